# 自动更新Excel文件 - daily refresh of "剩余" (E) / "开始时间" (F) columns.
#
# Business rule (reverse engineered from the commit): every day the sheet is
# refreshed, each data row's remaining-day counter (E) is recomputed as
#     E = D - (today - F)
# where D is the total day count and F is the last "start date" (stored as a
# literal YYYYMMDD integer, not an Excel date serial). Once the counter would
# hit zero or below, the row is considered "restocked" today, so its start
# date F is reset to today and the counter goes back up to the full D value.
#
# Reference date for this run is 2025-10-20 (per the commit message / diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = 20251020

# row -> @(new E value, new F value)
$updates = @{
    2  = @(14, 20251020)
    3  = @(14, 20251020)
    4  = @(14, 20251020)
    5  = @(6,  20251016)
    6  = @(14, 20251020)
    7  = @(6,  20251016)
    8  = @(14, 20251020)
    9  = @(6,  20251016)
    10 = @(7,  20251020)
    11 = @(14, 20251020)
    12 = @(6,  20251016)
    13 = @(14, 20251020)
    14 = @(14, 20251020)
    15 = @(14, 20251020)
    16 = @(10, 20251020)
    17 = @(6,  20251016)
    18 = @(9,  20251019)
    19 = @(9,  20251019)
    20 = @(9,  20251019)
    21 = @(9,  20251019)
    22 = @(6,  20251016)
    23 = @(6,  20251016)
    24 = @(6,  20251016)
    25 = @(6,  20251016)
    26 = @(6,  20251016)
    27 = @(1,  20251014)
    28 = @(9,  20251019)
    29 = @(9,  20251019)
    30 = @(9,  20251019)
    31 = @(9,  20251019)
    32 = @(9,  20251019)
    33 = @(9,  20251019)
    34 = @(9,  20251019)
    35 = @(9,  20251019)
    37 = @(9,  20251019)
    38 = @(9,  20251019)
    39 = @(9,  20251019)
    40 = @(7,  20251020)
    41 = @(7,  20251020)
    42 = @(9,  20251019)
    43 = @(6,  20251016)
    44 = @(7,  20251020)
    45 = @(6,  20251016)
    46 = @(7,  20251020)
    47 = @(9,  20251019)
    48 = @(7,  20251020)
    49 = @(1,  20251014)
    50 = @(4,  20251014)
    51 = @(4,  20251014)
    52 = @(4,  20251014)
    53 = @(4,  20251014)
    54 = @(4,  20251014)
    55 = @(4,  20251014)
    56 = @(4,  20251014)
    57 = @(4,  20251014)
    58 = @(8,  20251018)
    59 = @(8,  20251018)
    60 = @(8,  20251018)
    61 = @(1,  20251014)
    62 = @(8,  20251018)
    63 = @(8,  20251018)
    64 = @(8,  20251018)
    65 = @(9,  20251019)
    66 = @(9,  20251019)
    67 = @(9,  20251019)
    68 = @(9,  20251019)
    69 = @(9,  20251019)
    70 = @(10, 20251020)
    71 = @(10, 20251020)
    72 = @(10, 20251020)
    73 = @(10, 20251020)
    74 = @(10, 20251020)
    75 = @(10, 20251020)
    76 = @(10, 20251020)
    77 = @(3,  20251013)
    78 = @(3,  20251013)
    79 = @(3,  20251013)
    80 = @(3,  20251013)
    81 = @(3,  20251013)
    82 = @(3,  20251013)
    83 = @(3,  20251013)
    84 = @(3,  20251013)
    85 = @(3,  20251013)
    86 = @(3,  20251013)
    87 = @(7,  20251020)
    88 = @(7,  20251020)
    89 = @(7,  20251020)
    90 = @(7,  20251020)
    91 = @(6,  20251016)
    92 = @(7,  20251020)
    93 = @(3,  20251013)
    94 = @(3,  20251016)
    95 = @(2,  20251012)
    96 = @(10, 20251020)
    97 = @(10, 20251020)
    98 = @(10, 20251020)
    99 = @(10, 20251020)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 5).Value = $vals[0]
    $ws.Cells.Item($row, 6).Value = $vals[1]
}

Write-Host "Updated $($updates.Count) rows (columns E/F) to reflect $today refresh."
